# edit.ps1 -- applies the HUAG-014 user-story addition + HUAG 013 label
# re-run split described by the supplied diff, against Peaceathome.docx.
#
# Strategy: locate the two target paragraphs robustly via Find, then use
# Range.InsertXML (collapsed at the paragraph start) to replace each whole
# paragraph's contents with the exact target WordprocessingML. InsertXML
# on a collapsed range swaps out the entire enclosing <w:p> for the <w:p>
# element(s) supplied in the payload, leaving neighbouring paragraphs
# untouched -- which lets us reproduce the diff's run/bookmark structure
# exactly instead of relying on text-only Find/Replace (which would merge
# same-formatted runs into one <w:r>).

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Change 1: "HUAG 013: " (single run) -> "HUAG-" + "013: " (two runs)
# ---------------------------------------------------------------------

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("HUAG 013: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the 'HUAG 013: ' paragraph"
}
$huagPara = $rng1.Paragraphs(1)
$huagRange = $huagPara.Range
$huagRange.Collapse(1)

$huag013Xml = '<w:p><w:pPr><w:pStyle w:val="paragraph"/><w:spacing w:after="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>HUAG-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">013: </w:t></w:r></w:p>'

$huagRange.InsertXML((New-PkgXml $huag013Xml))

# ---------------------------------------------------------------------
# Change 2: replace the empty "_GoBack" paragraph (which follows the
# HUAG 013 story's last line) with the new HUAG-014 user story block.
# ---------------------------------------------------------------------

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("para poder resolver problemas", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the anchor paragraph preceding the _GoBack paragraph"
}
$anchorPara = $rng2.Paragraphs(1)
# The _GoBack paragraph is two paragraphs after the anchor: one blank
# spacer paragraph, then the (also blank) bookmark-only paragraph.
$blankPara = $anchorPara.Next()
$goBackPara = $blankPara.Next()

$goBackRange = $goBackPara.Range
$goBackRange.Collapse(1)

$boldPPr = '<w:pPr><w:pStyle w:val="paragraph"/><w:spacing w:after="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>'
$boldRPr = '<w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'

$plainPPrShort = '<w:pPr><w:pStyle w:val="paragraph"/><w:spacing w:after="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>'
$plainPPrFull = '<w:pPr><w:pStyle w:val="paragraph"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>'
$plainRPr = '<w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'

$p1 = '<w:p>' + $boldPPr + '<w:r>' + $boldRPr + '<w:t>HUAG-014</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r>' + $boldRPr + '<w:t>:</w:t></w:r></w:p>'

$p2 = '<w:p>' + $plainPPrShort + '<w:r>' + $plainRPr + '<w:t>Como t&#233;cnico de soporte,</w:t></w:r></w:p>'

$p3 = '<w:p>' + $plainPPrShort + '<w:proofErr w:type="gramStart"/><w:r>' + $plainRPr + '<w:t>quiero</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>' + $plainRPr + '<w:t xml:space="preserve"> ver una lista de tickets pendientes y asignarme uno,</w:t></w:r></w:p>'

$p4 = '<w:p>' + $plainPPrFull + '<w:proofErr w:type="gramStart"/><w:r>' + $plainRPr + '<w:t>para</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>' + $plainRPr + '<w:t xml:space="preserve"> poder empezar a resolverlo de forma organizada.</w:t></w:r></w:p>'

$p5 = '<w:p>' + $plainPPrFull + '</w:p>'

$p6 = '<w:p>' + $plainPPrFull + '</w:p>'

$goBackBlockXml = $p1 + $p2 + $p3 + $p4 + $p5 + $p6

$goBackRange.InsertXML((New-PkgXml $goBackBlockXml))

Write-Output "done"
